$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds text values that look numeric (e.g. "1.005").
# Excel's COM layer auto-coerces plain Value assignments that parse as a
# number, which would corrupt the original text-typed cells. Force the
# range to Text format first so the assigned strings stay strings, then
# clear the (temporary) number-format override so the cell style reverts
# to its original (unstyled) state once the values are in place.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.327.04"
$ws.Range("E2").Value = "  -7.49%  "
$ws.Range("D3").Value = "1.677.50"
$ws.Range("E3").Value = "  -6.26%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "218.61"
$ws.Range("E5").Value = "  -5.38%  "
$ws.Range("D6").Value = "0.5076"
$ws.Range("E6").Value = "  -13.54%  "
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.2662"
$ws.Range("E8").Value = "  -3.50%  "
$ws.Range("E9").Value = "  -5.34%  "
$ws.Range("D10").Value = "0.06317"
$ws.Range("E10").Value = "  -5.92%  "
$ws.Range("D11").Value = "0.07370"
$ws.Range("E11").Value = "  -2.27%  "
$ws.Range("D12").Value = "1.682.16"
$ws.Range("E12").Value = "  -5.99%  "
$ws.Range("D13").Value = "4.530"
$ws.Range("E13").Value = "  -4.89%  "
$ws.Range("D14").Value = "0.5789"
$ws.Range("E14").Value = "  -4.47%  "
$ws.Range("D15").Value = "1.909.92"
$ws.Range("E15").Value = "  -5.97%  "
$ws.Range("D16").Value = "0.000008522"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "64.70"
$ws.Range("E17").Value = "  -14.11%  "
$ws.Range("D18").Value = "26.375.77"
$ws.Range("E18").Value = "  -7.19%  "
$ws.Range("D19").Value = "5.006"
$ws.Range("E19").Value = "  -7.35%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  -4.63%  "
$ws.Range("D22").Value = "186.22"
$ws.Range("E22").Value = "  -10.26%  "
$ws.Range("D23").Value = "6.238"
$ws.Range("E23").Value = "  -7.87%  "
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D26").Value = "7.484"
$ws.Range("E26").Value = "  -7.13%  "
$ws.Range("E27").Value = "  -6.92%  "
$ws.Range("D28").Value = "15.85"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "1.344"
$ws.Range("E29").Value = "  -4.52%  "
$ws.Range("D30").Value = "0.05825"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").Value = "3.513"
$ws.Range("E32").Value = "  -6.74%  "
$ws.Range("D33").Value = "3.509"
$ws.Range("E33").Value = "  -6.63%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "0.5946"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("D37").Value = "2.361"
$ws.Range("E37").Value = "  -5.64%  "
$ws.Range("D38").Value = "2.672"
$ws.Range("E38").Value = "  -0.88%  "

# Rows 39/40 swap contents: Maker <-> VeChain (with updated price/volume).
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.094.99"
$ws.Range("E40").Value = "  -4.29%  "

$ws.Range("D41").Value = "5.892"
$ws.Range("E41").Value = "  -6.40%  "
$ws.Range("D42").Value = "0.8605"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "99.95"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").Value = "1.835.82"
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("D46").Value = "0.00000000117"
$ws.Range("E46").Value = "  +7.19%  "
$ws.Range("D47").Value = "56.27"
$ws.Range("E47").Value = "  -6.06%  "
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "8.021"
$ws.Range("E49").Value = "  -4.48%  "
$ws.Range("D50").Value = "0.4313"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").Value = "0.05207"
$ws.Range("E51").Value = "  -4.02%  "

# Drop the temporary Text number-format override so the cells end up back
# at the workbook's original (default/general) style, now holding text
# values instead of being re-coerced to numbers.
$priceRange.ClearFormats()

Write-Output "Applied cryptos price/volume update"
